$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-6 (columns A-I), column J keeps the shared "train_dim2_1" text
$data = @(
    @(1, 5, 8, 3, 4, -2, -4, 45, 5),
    @(2, 6, 7, 5, 2, -1, -5, 56, 5),
    @(3, 6, 9, 1, 8, -5, -1, 12, 5),
    @(4, 5, 6, 2, 3, -3, -3, 34, 5),
    @(5, 8, 9, 4, 7, -4, -2, 23, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$wb.Save()
